# Apply updated "dSF" (column F) values for specific rows on Sheet1,
# per the data repull / push-all-data / mean-calculation commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of worksheet row number -> new value for column F (dSF)
$updates = @{
    3  = 1
    5  = 0
    10 = -2
    11 = -1
    15 = 0
    20 = 2
    27 = 3
    28 = -4
    31 = 1
    35 = 4
    36 = 0
    37 = -1
    38 = -1
    40 = -2
    42 = 2
    47 = 0
    52 = -1
    63 = 2
    64 = -1
    66 = 0
    77 = -1
    79 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
